# Experiment using universal quantifier in ms power point.
#
# Slide 34: the "Comparing two implementations of a function" slide used the
# wrong glyph (U+2A5D, "⩝") where a universal quantifier (U+2200, "∀") was
# intended. Fix the glyph in place, keeping the Cambria Math run formatting.

$p = $ppt.ActivePresentation

$forall = [string][char]0x2200

function Set-CharRange($shape, [string]$needle, [string]$replacement) {
    $tr = $shape.TextFrame.TextRange
    $text = $tr.Text
    $idx = $text.IndexOf($needle)
    if ($idx -lt 0) {
        throw "needle not found: $needle"
    }
    $chars = $tr.Characters($idx + 1, $needle.Length)
    $chars.Text = $replacement
}

# --- Slide 34: "⩝ " -> "∀ " -------------------------------------------------
$s34 = $p.Slides.Item(34)
$shp34 = $s34.Shapes.Item(2)
$oldGlyph = [string][char]0x2A5D
Set-CharRange $shp34 ($oldGlyph + " ") ($forall + " ")

# --- Slide 35: merge "... E" + " " + "that transforms " runs ---------------
$s35 = $p.Slides.Item(35)
$shp35 = $s35.Shapes.Item(2)
Set-CharRange $shp35 " that transforms " " that transforms "

# --- Slide 41: merge "Could " + "be of " + "value for ..." runs ------------
$s41 = $p.Slides.Item(41)
$shp41 = $s41.Shapes.Item(2)
Set-CharRange $shp41 `
    "Could be of value for all kinds of mocking tools and Monte Carlo studies." `
    "Could be of value for all kinds of mocking tools and Monte Carlo studies."

# --- Slide 41: merge " formulae, commuting diagrams " + ". . .  " runs -----
Set-CharRange $shp41 " formulae, commuting diagrams . . .  " " formulae, commuting diagrams . . .  "

# --- Slide 42: merge "Transforms testing into " + "development" runs ------
$s42 = $p.Slides.Item(42)
$shp42 = $s42.Shapes.Item(2)
Set-CharRange $shp42 "Transforms testing into development" "Transforms testing into development"
